$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# URL value (row 2, column B)
$ws.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/CodeSystem/split-method"

# Version value (row 3, column B)
$ws.Range("B3").Value = "8.0.0"

# Date value (row 8, column B)
$ws.Range("B8").Value = "2022-11-10T16:00:46+00:00"

# Publisher value (row 9, column B)
$ws.Range("B9").Value = "LinuxForHealth Team"
